{"js": "// Restore the \"[DATE]\" placeholder in the Front Matter \"Status\" paragraph,\n// replacing the literal date \"August 7, 2023\" (the leading space before the\n// date is preserved, matching the author's original run boundary).\nconst body = context.document.body;\n\nconst dateResults = body.search(\" August 7, 2023\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error('Could not find \" August 7, 2023\" in the document body.');\n}\n\n// Replace the whole run's text (\" August 7, 2023\") with \" [DATE]\" in a single\n// call so the edit stays confined to that one run (no merging with the\n// neighboring runs before/after it).\nconst dateRun = dateResults.items[0];\ndateRun.insertText(\" [DATE]\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the inserted \"[DATE]\" text and nudge (then restore) one of its\n// formatting properties. This forces the host to materialize \"[DATE]\" as its\n// own run \u2014 distinct from the leading space run \u2014 mirroring the target\n// structure (two runs: \" \" and \"[DATE]\") instead of one combined \" [DATE]\" run.\nconst placeholderResults = body.search(\"[DATE]\", { matchCase: true });\nplaceholderResults.load(\"items\");\nawait context.sync();\n\nconst placeholderRun = placeholderResults.items[0];\nplaceholderRun.font.color = \"#FF0000\";\nawait context.sync();\n\nplaceholderRun.font.color = \"#000000\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Text = \"August 7, 2023\"\n$find.Replacement.Text = \"[DATE]\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
